# Rename the three worksheets ("DL" -> "CostStructure1", "JB" -> "CostStructure2",
# "MG" -> "CostStructure3") and update the saved selection (active cell) on the
# first and third sheets.

$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item("DL").Name = "CostStructure1"
$wb.Worksheets.Item("JB").Name = "CostStructure2"
$wb.Worksheets.Item("MG").Name = "CostStructure3"

$ws1 = $wb.Worksheets.Item("CostStructure1")
$ws1.Activate()
$ws1.Range("B25").Select()

$ws3 = $wb.Worksheets.Item("CostStructure3")
$ws3.Activate()
$ws3.Range("E15").Select()

$ws1.Activate()
